# feat: add 2022-Q3 data
#
# The workbook currently has 3 sheets: "总计" (summary), "2022-Q2", "2022-Q1".
# This script inserts a new "2022-Q3" sheet (a copy of the "2022-Q2" sheet,
# content replaced) right after "总计" / before "2022-Q2", and adds a
# corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

function Set-TextCell($rng, [string]$val) {
    # Force the cell to store $val as literal text (avoids Excel's
    # automatic "003308" -> 3308 / "4.20" -> 4.2 numeric coercion), then
    # drop the temporary number-format so the cell is left with the
    # default (no explicit) style, matching the rest of the sheet.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying "2022-Q2" (sheet 2) so it
#    inherits all styles/column widths/etc, placed right before it.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The copied sheet has 4 data rows (rows 2-5); the new data only has 2
# (rows 2-3), so drop the extra two rows.
$q3.Range("A4:H5").EntireRow.Delete()

# Row 2
$q3.Range("A2").Value = 0
Set-TextCell $q3.Range("B2") "003308"
Set-TextCell $q3.Range("C2") "中信建投睿利灵活配置混合A"
Set-TextCell $q3.Range("D2") "0.07"
Set-TextCell $q3.Range("E2") "93.78"
Set-TextCell $q3.Range("F2") "4.20"
Set-TextCell $q3.Range("G2") "0.0029"
$q3.Range("H2").Value = 5

# Row 3
$q3.Range("A3").Value = 1
Set-TextCell $q3.Range("B3") "004635"
Set-TextCell $q3.Range("C3") "中信建投睿利灵活配置混合C"
Set-TextCell $q3.Range("D3") "0.03"
Set-TextCell $q3.Range("E3") "93.78"
Set-TextCell $q3.Range("F3") "4.20"
Set-TextCell $q3.Range("G3") "0.0013"
$q3.Range("H3").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3,
#    pushing the existing 2022-Q2 / 2022-Q1 rows down, and renumber the
#    index column (A).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# Copy the index-column style (border/bold/center) from the row below
# onto the newly inserted A2 cell.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ------------------------------------------------------------------
# 3. Restore "2022-Q1" (now the 4th tab) as the active sheet, matching
#    the original workbook's active-tab state (unaffected by this edit).
# ------------------------------------------------------------------
$wb.Worksheets.Item(4).Activate()
